$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 170.71428
$ws.Range("I11").Value = 170.71428
$ws.Range("K11").Value = 170.71428
$ws.Range("M11").Value = -30.71428

$ws.Range("H31").Value = 125.5
$ws.Range("I31").Value = 125.5
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 376.5
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -146.5
$ws.Range("N31").ClearContents()

$ws.Range("H32").Value = 1167
$ws.Range("J32").Value = 1501
$ws.Range("L32").Value = 1501
$ws.Range("N32").Value = -2153

$ws.Range("H38").Value = 469
$ws.Range("I38").Value = 357.6
$ws.Range("J38").Value = 747.5
$ws.Range("K38").Value = 1072.8
$ws.Range("L38").Value = 2242.5
$ws.Range("M38").Value = -700.8000000000002
$ws.Range("N38").Value = -2986.5

$ws.Range("H39").Value = 110.166664
$ws.Range("I39").Value = 74.2
$ws.Range("J39").Value = 290
$ws.Range("K39").Value = 222.6
$ws.Range("L39").Value = 870
$ws.Range("M39").Value = 73.39999999999998
$ws.Range("N39").Value = -1462

$ws.Range("H111").Value = 480.625
$ws.Range("I111").Value = 335
$ws.Range("K111").Value = 1005
$ws.Range("M111").Value = 2062

$ws.Range("H115").Value = 1965
$ws.Range("I115").Value = 1965
$ws.Range("K115").Value = 5895
$ws.Range("M115").Value = -4328

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 537.75
$ws.Range("I2").Value = 537.75
$ws.Range("K2").Value = 537.75
$ws.Range("M2").Value = -424.75

$ws.Range("H6").Value = 650
$ws.Range("I6").Value = 600
$ws.Range("J6").Value = 700
$ws.Range("K6").Value = 600
$ws.Range("L6").Value = 700
$ws.Range("M6").Value = -427
$ws.Range("N6").Value = -1046

$ws.Range("H8").Value = 500
$ws.Range("I8").Value = 500
$ws.Range("K8").Value = 500
$ws.Range("M8").Value = -356

$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()

$ws.Range("H74").Value = 1830
$ws.Range("I74").Value = 1762.2
$ws.Range("K74").Value = 1762.2
$ws.Range("M74").Value = -888.2

$ws.Range("H77").Value = 1830
$ws.Range("I77").Value = 1762.2
$ws.Range("K77").Value = 8811
$ws.Range("M77").Value = -4443

$ws.Range("H116").Value = 537.75
$ws.Range("I116").Value = 537.75
$ws.Range("K116").Value = 537.75
$ws.Range("M116").Value = 1756.25

$ws.Range("H135").Value = 44956.332
$ws.Range("J135").Value = 44956.332
$ws.Range("L135").Value = 44956.332
$ws.Range("N135").Value = -55096.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 537.75
$ws.Range("I3").Value = 537.75
$ws.Range("K3").Value = 537.75
$ws.Range("M3").Value = -423.75

$ws.Range("H107").Value = 1669.7097
$ws.Range("I107").Value = 1493.2084
$ws.Range("K107").Value = 1493.2084
$ws.Range("M107").Value = 426.7916

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 62.857143
$ws.Range("I7").Value = 58.6
$ws.Range("J7").Value = 73.5
$ws.Range("K7").Value = 58.6
$ws.Range("L7").Value = 73.5
$ws.Range("M7").Value = 54.4
$ws.Range("N7").Value = -299.5

$ws.Range("H16").Value = 812
$ws.Range("I16").Value = 797.75
$ws.Range("J16").Value = 869
$ws.Range("K16").Value = 797.75
$ws.Range("L16").Value = 869
$ws.Range("M16").Value = -510.75
$ws.Range("N16").Value = -1443

$ws.Range("H22").Value = 706.17645
$ws.Range("I22").Value = 713.8
$ws.Range("J22").Value = 649
$ws.Range("K22").Value = 713.8
$ws.Range("L22").Value = 649
$ws.Range("M22").Value = -363.8
$ws.Range("N22").Value = -1349

$ws.Range("H32").Value = 2575
$ws.Range("I32").Value = 2575
$ws.Range("K32").Value = 2575
$ws.Range("M32").Value = -2259

$ws.Range("H113").Value = 812
$ws.Range("I113").Value = 797.75
$ws.Range("J113").Value = 869
$ws.Range("K113").Value = 797.75
$ws.Range("L113").Value = 869
$ws.Range("M113").Value = 1372.25
$ws.Range("N113").Value = -5209

$ws.Range("H133").Value = 74997
$ws.Range("J133").Value = 74997
$ws.Range("L133").Value = 74997
$ws.Range("N133").Value = -80057

$ws.Range("H134").Value = 1135.8889
$ws.Range("I134").Value = 1135.8889
$ws.Range("K134").Value = 3407.6667
$ws.Range("M134").Value = -872.6666999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2390.432
$ws.Range("I4").Value = 1896.0358
$ws.Range("J4").Value = 3255.625
$ws.Range("K4").Value = 5688.107400000001
$ws.Range("L4").Value = 9766.875
$ws.Range("M4").Value = -5576.107400000001
$ws.Range("N4").Value = -9990.875

$ws.Range("H44").Value = 1134.4445
$ws.Range("I44").Value = 617.5
$ws.Range("J44").Value = 1548
$ws.Range("K44").Value = 1852.5
$ws.Range("L44").Value = 4644
$ws.Range("M44").Value = -1454.5
$ws.Range("N44").Value = -5440

$ws.Range("H137").Value = 6749.5
$ws.Range("J137").Value = 6749.5
$ws.Range("L137").Value = 20248.5
$ws.Range("N137").Value = -30448.5

$ws.Range("H139").Value = 1687.5
$ws.Range("I139").Value = 1687.5
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 5062.5
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = 77.5
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 16333.333
$ws.Range("I43").Value = 7000
$ws.Range("J43").Value = 35000
$ws.Range("K43").Value = 7000
$ws.Range("L43").Value = 35000
$ws.Range("M43").Value = -6849
$ws.Range("N43").Value = -35302

$ws.Range("H47").Value = 30000
$ws.Range("J47").Value = 30000
$ws.Range("L47").Value = 30000
$ws.Range("N47").Value = -31136

$ws.Range("H55").Value = 33600
$ws.Range("J55").Value = 33600
$ws.Range("L55").Value = 33600
$ws.Range("N55").Value = -34254

$ws.Range("H102").Value = 59041.855
$ws.Range("I102").Value = 68615.5
$ws.Range("J102").Value = 1600
$ws.Range("K102").Value = 68615.5
$ws.Range("L102").Value = 1600
$ws.Range("M102").Value = -66993.5
$ws.Range("N102").Value = -4844

$ws.Range("H122").Value = 4655
$ws.Range("I122").Value = 4000
$ws.Range("J122").Value = 4873.3335
$ws.Range("K122").Value = 12000
$ws.Range("L122").Value = 14620.0005
$ws.Range("M122").Value = -9550
$ws.Range("N122").Value = -19520.0005

$ws.Range("H123").Value = 79966.664
$ws.Range("J123").Value = 79966.664
$ws.Range("L123").Value = 79966.664
$ws.Range("N123").Value = -84866.664

$ws.Range("H126").Value = 1480
$ws.Range("I126").Value = 1480
$ws.Range("K126").Value = 4440
$ws.Range("M126").Value = -1970

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2500
$ws.Range("I2").Value = 2500
$ws.Range("K2").Value = 2500
$ws.Range("M2").Value = -2388

$ws.Range("H13").Value = 7833.3335
$ws.Range("J13").Value = 7833.3335
$ws.Range("L13").Value = 7833.3335
$ws.Range("N13").Value = -8113.3335

$ws.Range("H16").Value = 3566.8333
$ws.Range("I16").Value = 3566.8333
$ws.Range("K16").Value = 3566.8333
$ws.Range("M16").Value = -3396.8333

$ws.Range("H22").Value = 2485.5
$ws.Range("I22").Value = 1220
$ws.Range("K22").Value = 1220
$ws.Range("M22").Value = -925

$ws.Range("H27").Value = 2485.5
$ws.Range("I27").Value = 1220
$ws.Range("K27").Value = 1220
$ws.Range("M27").Value = -1113

$ws.Range("H46").Value = 226377
$ws.Range("I46").Value = 500600
$ws.Range("J46").Value = 6998.6
$ws.Range("K46").Value = 500600
$ws.Range("L46").Value = 6998.6
$ws.Range("M46").Value = -500412
$ws.Range("N46").Value = -7374.6

$ws.Range("H68").Value = 2000
$ws.Range("I68").Value = 2000
$ws.Range("K68").Value = 2000
$ws.Range("M68").Value = -1251

$ws.Range("H71").Value = 2000
$ws.Range("I71").Value = 2000
$ws.Range("K71").Value = 10000
$ws.Range("M71").Value = -6256

$ws.Range("H133").Value = 100000
$ws.Range("J133").Value = 100000
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -105060

$ws.Range("H136").Value = 4697
$ws.Range("I136").Value = 4682.3335
$ws.Range("K136").Value = 14047.0005
$ws.Range("M136").Value = -11497.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H107").Value = 821.7143
$ws.Range("I107").Value = 360.1
$ws.Range("J107").Value = 1975.75
$ws.Range("K107").Value = 1080.3
$ws.Range("L107").Value = 5927.25
$ws.Range("M107").Value = 839.6999999999998
$ws.Range("N107").Value = -9767.25

$ws.Range("H132").Value = 1280
$ws.Range("I132").Value = 1280
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1310
$ws.Range("N132").ClearContents()
